# "Running all the test cases" — set the Runmode column (C2:C25) on the
# "Test Cases" sheet to "Y" for every test case, and select that range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

$range = $ws.Range("C2:C25")
$range.Value = "Y"
$range.Select()
